$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '26.997.46'
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = '  -1.41%  '
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.820.47'
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = '  -1.06%  '
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = '  -0.35%  '
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '310.21'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = '  -1.53%  '
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = '  -0.23%  '
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4632'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = '  -2.41%  '
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3625'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = '  -2.12%  '
$ws.Cells.Item(8, 5).Style = "Normal"
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = '  -2.41%  '
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.8647'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = '  -2.34%  '
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '19.81'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = '  -3.47%  '
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '1.857.62'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = '  -0.50%  '
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.07592'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = '  +2.98%  '
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '93.10'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = '  -0.17%  '
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '5.325'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = '  -2.98%  '
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '6.445'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = '  -2.06%  '
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '1.010'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = '  -0.40%  '
$ws.Cells.Item(17, 5).Style = "Normal"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.000008629'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = '  -2.42%  '
$ws.Cells.Item(18, 5).Style = "Normal"
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = '  -0.24%  '
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '27.311.85'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = '  -0.33%  '
$ws.Cells.Item(20, 5).Style = "Normal"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '14.44'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = '  -2.69%  '
$ws.Cells.Item(21, 5).Style = "Normal"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.153'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = '  -3.83%  '
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '10.56'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = '  -1.38%  '
$ws.Cells.Item(23, 5).Style = "Normal"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.083.35'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = '  +0.11%  '
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '151.21'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = '  -0.82%  '
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '1.866'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = '  -2.42%  '
$ws.Cells.Item(26, 5).Style = "Normal"
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = '  -2.33%  '
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.097'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = '  -3.35%  '
$ws.Cells.Item(28, 5).Style = "Normal"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '115.47'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = '  -2.14%  '
$ws.Cells.Item(29, 5).Style = "Normal"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '5.067'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = '  -3.71%  '
$ws.Cells.Item(30, 5).Style = "Normal"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.08887'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = '  -0.89%  '
$ws.Cells.Item(31, 5).Style = "Normal"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '2.958'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = '  +0.69%  '
$ws.Cells.Item(32, 5).Style = "Normal"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.7257'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = '  -4.56%  '
$ws.Cells.Item(33, 5).Style = "Normal"
$ws.Cells.Item(34, 5).NumberFormat = "@"
$ws.Cells.Item(34, 5).Value = '  -3.94%  '
$ws.Cells.Item(34, 5).Style = "Normal"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '4.415'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = '  -3.19%  '
$ws.Cells.Item(35, 5).Style = "Normal"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.012'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = '  -0.09%  '
$ws.Cells.Item(36, 5).Style = "Normal"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.514'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = '  +5.80%  '
$ws.Cells.Item(37, 5).Style = "Normal"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '1.071'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = '  -3.31%  '
$ws.Cells.Item(38, 5).Style = "Normal"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.05250'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = '  -2.37%  '
$ws.Cells.Item(39, 5).Style = "Normal"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.01910'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = '  -2.58%  '
$ws.Cells.Item(40, 5).Style = "Normal"
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = '  -2.67%  '
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '7.088'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = '  -2.75%  '
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.5194'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = '  -3.05%  '
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.1627'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = '  -2.27%  '
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '8.207'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = '  -4.04%  '
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.4849'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = '  -2.51%  '
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.012'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = '  -0.24%  '
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '10.15'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = '  -3.64%  '
$ws.Cells.Item(48, 5).Style = "Normal"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '103.27'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = '  -1.82%  '
$ws.Cells.Item(49, 5).Style = "Normal"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.634'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = '  -2.71%  '
$ws.Cells.Item(50, 5).Style = "Normal"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.06240'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = '  -1.22%  '
$ws.Cells.Item(51, 5).Style = "Normal"
